$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.05
$ws.Range("O2").Value = 1.33
$ws.Range("M3").Value = 1.05
$ws.Range("O3").Value = 1.37
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.67
$ws.Range("M4").Value = 1.05
$ws.Range("O4").Value = 1.41
$ws.Range("P4").Value = 2.7
$ws.Range("G5").Value = 2.5
$ws.Range("H5").Value = 2.7
$ws.Range("G6").Value = 1.76
$ws.Range("G7").Value = 2.3
$ws.Range("I7").Value = 2.5
$ws.Range("K7").Value = 2.25
$ws.Range("Q7").Value = 1.75
$ws.Range("R7").Value = 2.05
$ws.Range("U7").Value = 1.67
$ws.Range("V7").Value = 2.1
$ws.Range("W7").Value = 9.5
$ws.Range("AE7").Value = 13
$ws.Range("AG7").Value = 151
$ws.Range("AH7").Value = 10
$ws.Range("AZ7").Value = 41
$ws.Range("BA7").Value = 51
$ws.Range("BB7").Value = 126
$ws.Range("BD7").Value = 151
$ws.Range("I8").Value = 1.77
$ws.Range("O10").Value = 1.11
$ws.Range("G11").Value = 1.71
$ws.Range("Q11").Value = 1.57
$ws.Range("G12").Value = 2.2
$ws.Range("Q12").Value = 1.63
$ws.Range("I13").Value = 1.79
$ws.Range("Q13").Value = 1.37
$ws.Range("R13").Value = 2.87
$ws.Range("G17").Value = 1.5
$ws.Range("Q17").Value = 1.82
$ws.Range("R17").Value = 1.92
$ws.Range("I18").Value = 2.87
$ws.Range("G20").Value = 1.69
$ws.Range("U25").Value = 1.77
$ws.Range("V25").Value = 1.92
$ws.Range("U26").Value = 1.58
$ws.Range("O27").Value = 1.36
$ws.Range("P27").Value = 3
$ws.Range("Q27").Value = 2.2
$ws.Range("R27").Value = 1.65
$ws.Range("V27").Value = 1.69
$ws.Range("U28").Value = 1.69
$ws.Range("M30").Value = 1.05
$ws.Range("N30").Value = 11
$ws.Range("Q30").Value = 1.98
$ws.Range("R30").Value = 1.88
$ws.Range("G31").Value = 2.35
$ws.Range("I31").Value = 2.75
$ws.Range("K31").Value = 2.4
$ws.Range("N31").Value = 17
$ws.Range("Q31").Value = 1.53
$ws.Range("R31").Value = 2.4
$ws.Range("AE31").Value = 11
$ws.Range("AJ31").Value = 11
$ws.Range("AP31").Value = 17
$ws.Range("BC31").Value = 301
$ws.Range("O32").Value = 1.36
$ws.Range("P32").Value = 3
$ws.Range("Q33").Value = 1.7
$ws.Range("R33").Value = 2.1
$ws.Range("M34").Value = 1.07
$ws.Range("N34").Value = 9
$ws.Range("Q35").Value = 1.92
$ws.Range("R35").Value = 1.82
